$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = '''43.909.88'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +0.46%  '

$c = $ws.Range("D3")
$c.Value = '''2.294.32'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -0.66%  '

$c = $ws.Range("D4")
$c.Value = '''0.999'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.26%  '

$c = $ws.Range("D5")
$c.Value = '''107.90'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +11.21%  '

$c = $ws.Range("D6")
$c.Value = '''271.74'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.36%  '

$c = $ws.Range("D7")
$c.Value = '''0.626'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.33%  '

$ws.Range("E8").Value = '  +0.16%  '

$c = $ws.Range("D9")
$c.Value = '''0.612'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -2.08%  '

$c = $ws.Range("D10")
$c.Value = '''46.38'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +2.55%  '

$c = $ws.Range("D11")
$c.Value = '''0.0937'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -1.59%  '

$c = $ws.Range("D12")
$c.Value = '''8.28'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +3.33%  '

$c = $ws.Range("D13")
$c.Value = '''0.108'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +1.94%  '

$c = $ws.Range("D14")
$c.Value = '''15.60'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +0.53%  '

$c = $ws.Range("D15")
$c.Value = '''2.637.08'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -0.41%  '

$c = $ws.Range("D16")
$c.Value = '''0.856'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -2.06%  '

$c = $ws.Range("D17")
$c.Value = '''2.291.93'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -1.04%  '

$c = $ws.Range("D18")
$c.Value = '''43.786.18'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.21%  '

$ws.Range("E19").Value = '  +0.33%  '

$c = $ws.Range("D20")
$c.Value = '''6.30'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -1.93%  '

$c = $ws.Range("D21")
$c.Value = '''72.18'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -1.77%  '

$c = $ws.Range("D22")
$c.Value = '''2.50'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +9.83%  '

$c = $ws.Range("D23")
$c.Value = '''233.67'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -2.47%  '

$c = $ws.Range("D24")
$c.Value = '''2.94'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +15.95%  '

$c = $ws.Range("D25")
$c.Value = '''9.29'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -1.65%  '

$ws.Range("E26").Value = '  +0.05%  '

$c = $ws.Range("D27")
$c.Value = '''11.32'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.24%  '

$ws.Range("B28").Value = 'InjectiveProtocol'
$ws.Range("C28").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range("D28")
$c.Value = '''40.89'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +6.85%  '

$ws.Range("B29").Value = 'WEMIXToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range("D29")
$c.Value = '''3.45'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -1.00%  '

$c = $ws.Range("D30")
$c.Value = '''2.28'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.40%  '

$c = $ws.Range("D31")
$c.Value = '''177.84'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +1.49%  '

$c = $ws.Range("D32")
$c.Value = '''21.87'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -2.67%  '

$c = $ws.Range("D33")
$c.Value = '''0.0909'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -0.68%  '

$c = $ws.Range("D34")
$c.Value = '''5.55'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +1.16%  '

$c = $ws.Range("D35")
$c.Value = '''4.91'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +10.31%  '

$ws.Range("E36").Value = '  -0.34%  '

$c = $ws.Range("D37")
$c.Value = '''0.113'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +3.19%  '

$c = $ws.Range("D38")
$c.Value = '''0.0361'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -1.35%  '

$c = $ws.Range("D39")
$c.Value = '''3.60'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +6.73%  '

$c = $ws.Range("D40")
$c.Value = '''0.237'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -3.22%  '

$c = $ws.Range("D41")
$c.Value = '''2.34'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -2.66%  '

$c = $ws.Range("D42")
$c.Value = '''1.38'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -2.83%  '

$c = $ws.Range("D43")
$c.Value = '''66.11'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +5.21%  '

$c = $ws.Range("D44")
$c.Value = '''12.26'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.73%  '

$c = $ws.Range("D45")
$c.Value = '''5.47'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +2.39%  '

$c = $ws.Range("D46")
$c.Value = '''8.80'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -4.22%  '

$ws.Range("E47").Value = '  -1.92%  '

$c = $ws.Range("D48")
$c.Value = '''1.24'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +2.61%  '

$c = $ws.Range("D49")
$c.Value = '''99.61'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -0.90%  '

$ws.Range("E50").Value = '  +11.73%  '

$c = $ws.Range("D51")
$c.Value = '''0.440'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +4.98%  '
